# correction lien vers profinfo
# Swap the student (No de DA + Nom de l'etudiant) between row 2 and row 11
# on the "rencontre-fin" sheet, leaving the date/time column (A) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rencontre-fin")

# Capture current values/formulas for B2:C2 and B11:C11 before swapping
$b2Formula = $ws.Range("B2").Formula
$c2Value   = $ws.Range("C2").Value2

$b11Formula = $ws.Range("B11").Formula
$c11Value   = $ws.Range("C11").Value2

# Write the swapped values back
$ws.Range("B2").Formula  = $b11Formula
$ws.Range("C2").Value2   = $c11Value

$ws.Range("B11").Formula = $b2Formula
$ws.Range("C11").Value2  = $c2Value

# Update the active selection to match the edited cells
$ws.Activate()
$ws.Range("B11:C11").Select()
